$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45760
$ws.Range("B2").Value = 2.771428571428572

$ws.Range("A3").Value = 45767
$ws.Range("B3").Value = 3.6

$ws.Range("A4").Value = 45774
$ws.Range("B4").Value = 0.1285714285714286

$ws.Range("A5").Value = 45781
$ws.Range("B5").Value = 6.471428571428571
